# Horarios actualizados Linea 141 - 1163
# Applies the latest scrape snapshot to each schedule sheet:
#   - refreshes the "Ultima actualizacion" / "Total filas" header cells
#   - updates rows whose data changed between scrapes (timestamp/line/minutes)
#   - appends newly scraped rows at the bottom of each table
$wb = $excel.ActiveWorkbook

# ----- Sheet: LP1912 -----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:59:02"
$ws.Cells.Item(3, 1).Value = "Total filas: 234"
$ws.Cells.Item(22, 3).Value = "215C_EL PATO"
$ws.Cells.Item(23, 3).Value = "14_ABASTO"
$ws.Cells.Item(62, 1).Value = "08:04:39"
$ws.Cells.Item(62, 3).Value = "14_ABASTO"
$ws.Cells.Item(62, 4).Value = 39
$ws.Cells.Item(63, 1).Value = "08:41:16"
$ws.Cells.Item(63, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(63, 4).Value = 2
$ws.Cells.Item(76, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(77, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(82, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(84, 3).Value = "17_ROMERO"
$ws.Cells.Item(98, 1).Value = "08:41:16"
$ws.Cells.Item(98, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(98, 4).Value = 90
$ws.Cells.Item(99, 1).Value = "09:01:18"
$ws.Cells.Item(99, 3).Value = "10_OLMOS"
$ws.Cells.Item(99, 4).Value = 70
$ws.Cells.Item(118, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(119, 3).Value = "10_OLMOS"
$ws.Cells.Item(141, 1).Value = "11:54:47"
$ws.Cells.Item(141, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(141, 4).Value = 12
$ws.Cells.Item(142, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(143, 1).Value = "11:07:42"
$ws.Cells.Item(143, 3).Value = "14_ABASTO"
$ws.Cells.Item(143, 4).Value = 59
$ws.Cells.Item(155, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(156, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(169, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(170, 3).Value = "215D_EL PATO"
$ws.Cells.Item(186, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(187, 3).Value = "17_ROMERO"
$ws.Cells.Item(203, 1).Value = "14:59:02"
$ws.Cells.Item(203, 2).Value = "14:59"
$ws.Cells.Item(203, 3).Value = "215B_EL PATO"
$ws.Cells.Item(203, 4).Value = 0
$ws.Cells.Item(204, 1).Value = "14:59:02"
$ws.Cells.Item(204, 2).Value = "15:00"
$ws.Cells.Item(204, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(204, 4).Value = 1
$ws.Cells.Item(205, 1).Value = "14:59:02"
$ws.Cells.Item(205, 2).Value = "15:05"
$ws.Cells.Item(205, 3).Value = "10_OLMOS"
$ws.Cells.Item(205, 4).Value = 6
$ws.Cells.Item(206, 1).Value = "14:59:02"
$ws.Cells.Item(206, 2).Value = "15:05"
$ws.Cells.Item(206, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(206, 4).Value = 6
$ws.Cells.Item(207, 1).Value = "14:59:02"
$ws.Cells.Item(207, 2).Value = "15:07"
$ws.Cells.Item(207, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(207, 4).Value = 8
$ws.Cells.Item(208, 1).Value = "14:59:02"
$ws.Cells.Item(208, 2).Value = "15:10"
$ws.Cells.Item(208, 3).Value = "17_ROMERO"
$ws.Cells.Item(208, 4).Value = 11
$ws.Cells.Item(209, 1).Value = "14:24:30"
$ws.Cells.Item(209, 2).Value = "15:13"
$ws.Cells.Item(209, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(209, 4).Value = 49
$ws.Cells.Item(210, 1).Value = "14:59:02"
$ws.Cells.Item(210, 2).Value = "15:14"
$ws.Cells.Item(210, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(210, 4).Value = 15
$ws.Cells.Item(211, 1).Value = "14:59:02"
$ws.Cells.Item(211, 2).Value = "15:17"
$ws.Cells.Item(211, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(211, 4).Value = 18
$ws.Cells.Item(212, 1).Value = "14:59:02"
$ws.Cells.Item(212, 2).Value = "15:20"
$ws.Cells.Item(212, 3).Value = "15_ABASTO"
$ws.Cells.Item(212, 4).Value = 21
$ws.Cells.Item(213, 1).Value = "14:59:02"
$ws.Cells.Item(213, 2).Value = "15:21"
$ws.Cells.Item(213, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(213, 4).Value = 22
$ws.Cells.Item(214, 1).Value = "13:45:48"
$ws.Cells.Item(214, 2).Value = "15:22"
$ws.Cells.Item(214, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(214, 4).Value = 97
$ws.Cells.Item(215, 1).Value = "14:59:02"
$ws.Cells.Item(215, 2).Value = "15:32"
$ws.Cells.Item(215, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(215, 4).Value = 33
$ws.Cells.Item(216, 1).Value = "13:45:48"
$ws.Cells.Item(216, 2).Value = "15:34"
$ws.Cells.Item(216, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(216, 4).Value = 109
$ws.Cells.Item(217, 1).Value = "14:59:02"
$ws.Cells.Item(217, 2).Value = "15:35"
$ws.Cells.Item(217, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(217, 4).Value = 36
$ws.Cells.Item(218, 1).Value = "14:59:02"
$ws.Cells.Item(218, 2).Value = "15:37"
$ws.Cells.Item(218, 3).Value = "10_OLMOS"
$ws.Cells.Item(218, 4).Value = 38
$ws.Cells.Item(219, 2).Value = "15:38"
$ws.Cells.Item(219, 3).Value = "215A_EL PATO"
$ws.Cells.Item(219, 4).Value = 74
$ws.Cells.Item(220, 2).Value = "15:38"
$ws.Cells.Item(220, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(220, 4).Value = 74
$ws.Cells.Item(221, 1).Value = "14:59:02"
$ws.Cells.Item(221, 2).Value = "15:39"
$ws.Cells.Item(221, 3).Value = "215A_EL PATO"
$ws.Cells.Item(221, 4).Value = 40
$ws.Cells.Item(222, 1).Value = "13:45:48"
$ws.Cells.Item(222, 2).Value = "15:42"
$ws.Cells.Item(222, 3).Value = "14_ABASTO"
$ws.Cells.Item(222, 4).Value = 117
$ws.Cells.Item(223, 2).Value = "15:45"
$ws.Cells.Item(223, 3).Value = "14_ABASTO"
$ws.Cells.Item(223, 4).Value = 81
$ws.Cells.Item(224, 1).Value = "14:24:30"
$ws.Cells.Item(224, 2).Value = "15:46"
$ws.Cells.Item(224, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(224, 4).Value = 82
$ws.Cells.Item(224, 5).Value = "LP1912"
$ws.Cells.Item(225, 1).Value = "14:59:02"
$ws.Cells.Item(225, 2).Value = "15:46"
$ws.Cells.Item(225, 3).Value = "14_ABASTO"
$ws.Cells.Item(225, 4).Value = 47
$ws.Cells.Item(225, 5).Value = "LP1912"
$ws.Cells.Item(226, 1).Value = "14:59:02"
$ws.Cells.Item(226, 2).Value = "15:47"
$ws.Cells.Item(226, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(226, 4).Value = 48
$ws.Cells.Item(226, 5).Value = "LP1912"
$ws.Cells.Item(227, 1).Value = "14:24:30"
$ws.Cells.Item(227, 2).Value = "15:53"
$ws.Cells.Item(227, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(227, 4).Value = 89
$ws.Cells.Item(227, 5).Value = "LP1912"
$ws.Cells.Item(228, 1).Value = "14:59:02"
$ws.Cells.Item(228, 2).Value = "15:54"
$ws.Cells.Item(228, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(228, 4).Value = 55
$ws.Cells.Item(228, 5).Value = "LP1912"
$ws.Cells.Item(229, 1).Value = "14:24:30"
$ws.Cells.Item(229, 2).Value = "15:56"
$ws.Cells.Item(229, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(229, 4).Value = 92
$ws.Cells.Item(229, 5).Value = "LP1912"
$ws.Cells.Item(230, 1).Value = "14:59:02"
$ws.Cells.Item(230, 2).Value = "15:56"
$ws.Cells.Item(230, 3).Value = "17_ROMERO"
$ws.Cells.Item(230, 4).Value = 57
$ws.Cells.Item(230, 5).Value = "LP1912"
$ws.Cells.Item(231, 1).Value = "14:59:02"
$ws.Cells.Item(231, 2).Value = "15:57"
$ws.Cells.Item(231, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(231, 4).Value = 58
$ws.Cells.Item(231, 5).Value = "LP1912"
$ws.Cells.Item(232, 1).Value = "14:59:02"
$ws.Cells.Item(232, 2).Value = "16:09"
$ws.Cells.Item(232, 3).Value = "14_ABASTO"
$ws.Cells.Item(232, 4).Value = 70
$ws.Cells.Item(232, 5).Value = "LP1912"
$ws.Cells.Item(233, 1).Value = "14:59:02"
$ws.Cells.Item(233, 2).Value = "16:15"
$ws.Cells.Item(233, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(233, 4).Value = 76
$ws.Cells.Item(233, 5).Value = "LP1912"
$ws.Cells.Item(234, 1).Value = "14:59:02"
$ws.Cells.Item(234, 2).Value = "16:20"
$ws.Cells.Item(234, 3).Value = "215C_EL PATO"
$ws.Cells.Item(234, 4).Value = 81
$ws.Cells.Item(234, 5).Value = "LP1912"
$ws.Cells.Item(235, 1).Value = "14:59:02"
$ws.Cells.Item(235, 2).Value = "16:21"
$ws.Cells.Item(235, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(235, 4).Value = 82
$ws.Cells.Item(235, 5).Value = "LP1912"
$ws.Cells.Item(236, 1).Value = "14:59:02"
$ws.Cells.Item(236, 2).Value = "16:30"
$ws.Cells.Item(236, 3).Value = "15_ABASTO"
$ws.Cells.Item(236, 4).Value = 91
$ws.Cells.Item(236, 5).Value = "LP1912"
$ws.Cells.Item(237, 1).Value = "14:59:02"
$ws.Cells.Item(237, 2).Value = "16:43"
$ws.Cells.Item(237, 3).Value = "225_GOMEZ"
$ws.Cells.Item(237, 4).Value = 104
$ws.Cells.Item(237, 5).Value = "LP1912"
$ws.Cells.Item(238, 1).Value = "14:59:02"
$ws.Cells.Item(238, 2).Value = "16:43"
$ws.Cells.Item(238, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(238, 4).Value = 104
$ws.Cells.Item(238, 5).Value = "LP1912"
$ws.Cells.Item(239, 1).Value = "14:59:02"
$ws.Cells.Item(239, 2).Value = "16:56"
$ws.Cells.Item(239, 3).Value = "17_179 Y 38"
$ws.Cells.Item(239, 4).Value = 117
$ws.Cells.Item(239, 5).Value = "LP1912"

# ----- Sheet: LP1912-215 -----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:59:02"
$ws.Cells.Item(3, 1).Value = "Total filas: 26"
$ws.Cells.Item(28, 1).Value = "14:59:02"
$ws.Cells.Item(28, 2).Value = "14:59"
$ws.Cells.Item(28, 3).Value = "215B_EL PATO"
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 2).Value = "15:38"
$ws.Cells.Item(29, 3).Value = "215A_EL PATO"
$ws.Cells.Item(29, 4).Value = 74
$ws.Cells.Item(30, 1).Value = "14:59:02"
$ws.Cells.Item(30, 2).Value = "15:39"
$ws.Cells.Item(30, 3).Value = "215A_EL PATO"
$ws.Cells.Item(30, 4).Value = 40
$ws.Cells.Item(30, 5).Value = "LP1912"
$ws.Cells.Item(31, 1).Value = "14:59:02"
$ws.Cells.Item(31, 2).Value = "16:20"
$ws.Cells.Item(31, 3).Value = "215C_EL PATO"
$ws.Cells.Item(31, 4).Value = 81
$ws.Cells.Item(31, 5).Value = "LP1912"

# ----- Sheet: 6203-6173 -----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:59:02"
$ws.Cells.Item(3, 1).Value = "Total filas: 32"
$ws.Cells.Item(35, 1).Value = "14:59:02"
$ws.Cells.Item(35, 4).Value = 35
$ws.Cells.Item(36, 1).Value = "14:59:02"
$ws.Cells.Item(36, 4).Value = 75
$ws.Cells.Item(37, 1).Value = "14:59:02"
$ws.Cells.Item(37, 2).Value = "16:53"
$ws.Cells.Item(37, 3).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(37, 4).Value = 114
$ws.Cells.Item(37, 5).Value = "L6173"

